# Applies the crew-list data edits described by the diff:
# ten rows in the single table each get a new name, nationality,
# date/place of birth, and ID code.

$d = $word.ActiveDocument

function Replace-Unique($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 0, $false, $replace, 1) | Out-Null
}

function Replace-InCell($table, $rowIndex, $colIndex, $find, $replace) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($colIndex)
    $range = $d.Range($cell.Range.Start, $cell.Range.End)
    $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 0, $false, $replace, 1) | Out-Null
}

# Row: Garza, Daumantas Sawney
Replace-Unique "Garza, Daumantas Sawney" "Meyrick, BertrandLovise"
Replace-Unique "ITA" "LCA"
Replace-Unique "1988-11-30, London" "1990-12-16, Nairobi"
Replace-Unique "sUPnZvnwYL" "Lr8mBbKmDX"

# Row: Rake, Batuhan Khodadad
Replace-Unique "Rake, Batuhan Khodadad" "Genadiev, AyselKamil"
Replace-Unique "PRT" "LVA"
Replace-Unique "1994-11-02, Riyadh" "1975-07-12, Basra"
Replace-Unique "ZNhrhsG6qh" "ABRRx5vU4R"

# Row: Evangelista, Veniamin Kyros
Replace-Unique "Evangelista, Veniamin Kyros" "Traylor, GovindaDiodotos"
Replace-Unique "GAB" "TLS"
Replace-Unique "1966-10-12, Yaounde" "1961-08-01, Rome"
Replace-Unique "MmdApg7tJA" "9WD3JqhjZ8"

# Row: Crespo, HonorataMarko
Replace-Unique "Crespo, HonorataMarko" "Hakim, KalinLalita"
Replace-Unique "MDA" "CXR"
Replace-Unique "1985-03-31, Dallas" "1970-01-04, Kanpur"
Replace-Unique "T6Dj6ky8KQ" "2HNrpQt8W5"

# Row: Furlan, IvkaWilla (table row 23 - nationality "MDG" is ambiguous doc-wide)
$table = $d.Tables.Item(1)
Replace-Unique "Furlan, IvkaWilla" "Post, Eun-JiIryna"
Replace-InCell $table 23 4 "MDG" "GUM"
Replace-Unique "1965-03-14, Palembang" "1968-07-14, Munich"
Replace-Unique "vxVbNYPs73" "UhehqHBRCf"

# Row: MacNevin, ApostolKanti
Replace-Unique "MacNevin, ApostolKanti" "Sanna, GuilhermeRamana"
Replace-Unique "MAR" "BMU"
Replace-Unique "1992-01-09, Philadelphia" "1995-10-12, Pyongyang"
Replace-Unique "H3nnfEZZaB" "vVwAtRS6LV"

# Row: Amador, DaudZinat
Replace-Unique "Amador, DaudZinat" "Aitken, AdelaisEsdras"
Replace-Unique "TJK" "WSM"
Replace-Unique "1965-06-12, Zhengzhou" "1988-10-17, Shanghai"
Replace-Unique "2E553AnFAJ" "sNneLaZsZd"

# Row: Christian, VilhjálmurTerje
Replace-Unique "Christian, VilhjálmurTerje" "Chiara, AniMaria"
Replace-Unique "FJI" "NGA"
Replace-Unique "1988-11-25, Fuzhou" "1989-11-14, Lahore"
Replace-Unique "uFqehB88At" "gGkrduyGB2"

# Row: Alessi, AuroraAmonet
Replace-Unique "Alessi, AuroraAmonet" "Wheatley, MelisaViraj"
Replace-Unique "MHL" "COG"
Replace-Unique "1958-08-10, Kuala Lumpur" "1959-06-09, Minsk"
Replace-Unique "BPFPFALjJS" "qpTEBJp8f2"

# Row: Gwerder, MileJuliana (table row 38 - rank "UXO" and nationality "MDG" are
# ambiguous doc-wide, so scope those two replacements to the specific cells)
Replace-Unique "Gwerder, MileJuliana" "Bridges, LauriLakshmana"
Replace-InCell $table 38 3 "UXO" "Client"
Replace-InCell $table 38 4 "MDG" "NGA"
Replace-Unique "2000-08-28, Caracas" "1986-04-08, Prague"
Replace-Unique "URgATaEsgN" "JrDrsZkesB"

Write-Output "Edit complete"
